# edit.ps1 - apply resume content updates per commit message / diff
$d = $word.ActiveDocument

function Replace-Text($oldText, $newText) {
    $ok = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                                   $true, 1, $false, $newText, 2)
    if (-not $ok) {
        Write-Output "WARNING: replace not found -> $oldText"
    }
}

function Insert-ParagraphAfter($anchorText, $newParaText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Output "WARNING: anchor not found -> $anchorText"
        return
    }
    $endPos = $rng.End
    $breakPoint = $d.Range($endPos, $endPos)
    $breakPoint.InsertParagraphAfter()
    $newPara = $d.Range($endPos + 1, $endPos + 1)
    $newPara.InsertAfter($newParaText)
}

# 1. Update years of experience in the summary paragraph
Replace-Text "Research & Data Professional with 21 years of experience" `
             "Research & Data Professional with 15+ years of experience"

# 2. Enhance FLEEM bullet (Progressive Change Campaign Committee)
Replace-Text "• Conceived, architected, and engineered FLEEM web application using Twilio API for thousands of simultaneous phone calls" `
             "• Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of calls using emulated predictive dialer for regulated political surveys"

# 3. Enhance Salsa Labs CRM bullet
Replace-Text "• Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system" `
             "• Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system used by tens of thousands of users simultaneously"

# 4. Enhance Salsa Labs mapping bullet
Replace-Text "• Integrated mapping and visualization tools for political campaign data analysis" `
             "• Integrated mapping and visualization tools for political campaign data analysis interfacing with Government and Activism APIs"

# 5. New bullet after "Collaborated with political strategists..." (end of Salsa Labs section)
Insert-ParagraphAfter "• Collaborated with political strategists to translate geospatial requirements into technical solutions" `
                       "• Handled billions of records with millions of columns in high-performance CRM system"

# 6. New bullet at end of Praxis Project section (before Lake Research Partners heading)
Insert-ParagraphAfter "• Managed technology infrastructure supporting community health initiatives across multiple countries" `
                       "• Architected and developed 25 Drupal sites to integrate with membership databases, activism CRMs and government agencies, under guidelines from Kellogg Foundation and Robert Wood Johnson Foundation"

# 7. New bullet at end of Lake Research Partners section (before Feldman Group heading)
Insert-ParagraphAfter "• Developed innovative approaches to visualizing demographic and market data for enhanced client understanding" `
                       "• Trained staff on building Python tooling for report generation and analysis"

# 8. New bullet at end of Feldman Group section (before Key Achievements heading)
Insert-ParagraphAfter "• Enhanced value of research deliverables through advanced analytical techniques using SPSS, OSCAR, PHP, and MySQL" `
                       "• Trained staff on PHP/MySQL for data analysis and reporting systems"
